# Comandos TypeScript.xlsx - add "Comandos Node" section (rows 20-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New section header: "Comandos Node" -----------------------------------
$ws.Range("A20").Value = "Comandos Node"

# npm init
$ws.Range("A21").Value = "npm init"
$ws.Range("B21").Value = "Iniciar un proyecto de node.js"

# npm install lite-server --save-dev (+ link to the project on GitHub)
$ws.Range("A22").Value = "npm install lite-server --save-dev"
$ws.Range("B22").Value = "Instalación de un servidor ligero"
$ws.Range("D22").Value = "https://github.com/johnpapa/lite-server"

# npm run dev
$ws.Range("A23").Value = "npm run dev"
$ws.Range("B23").Value = "Ejecuta el servidor configurado"

# npm i npm
$ws.Range("A24").Value = "npm i npm"
$ws.Range("B24").Value = "Instala las dependencias de un proyecto"

# npm install -g lite-server (+ link to the project on GitHub)
$ws.Range("A25").Value = "npm install -g lite-server"
$ws.Range("B25").Value = "Para instalar globalmente"
$ws.Range("D25").Value = "https://github.com/johnpapa/lite-server"

# --- Hyperlinks for the two new GitHub references --------------------------
$null = $ws.Hyperlinks.Add($ws.Range("D22"), "https://github.com/johnpapa/lite-server")
$ws.Range("D22").Style = "Hyperlink"

$null = $ws.Hyperlinks.Add($ws.Range("D25"), "https://github.com/johnpapa/lite-server")
$ws.Range("D25").Style = "Hyperlink"

# --- Page setup: portrait orientation --------------------------------------
$ws.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait

# --- View state: scroll position / active selection -------------------------
$excel.ActiveWindow.ScrollRow = 5
$ws.Range("A26").Select()

Write-Host "Applied Comandos Node section."
